# Update the "Förändrad" (last-changed) date in column C for every data
# row (rows 2-400) from Excel serial date 45190 (2023-09-21) to
# 45192 (2023-09-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("C2:C400")
$rng.Value2 = 45192
